# Weekly reports and Timesheets
# Weekly reports and weekly timesheets. Made some dates corrections on
# some weekly Timesheets
#
# This timesheet previously covered the week of 9/02-15/02; it is being
# rolled forward two weeks to 23/02-01/03.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# "Week of:" date (top of sheet)
$ws.Range("G8").Value = 41693

# Day-of-week column labels for the daily hours table
$ws.Range("A11").Value = "Sun 23/02"
$ws.Range("A12").Value = "Mon 24/02"
$ws.Range("A13").Value = "Tue 25/02"
$ws.Range("A14").Value = "Wed  26/02"
$ws.Range("A15").Value = "Thur  27/02"
$ws.Range("A16").Value = "Fri   28/02"
$ws.Range("A17").Value = "Sat  01/03"

# Employee signature date (end of the new week)
$ws.Range("D25").Value = 41700

# Supervisor signature date (unchanged, kept in sync with the form)
$ws.Range("D27").Value = 41728

# Update the on-screen selection / scroll position to match where the
# author was working when the sheet was last saved.
$ws.Activate() | Out-Null
$ws.Range("D25:E25").Select() | Out-Null
